# Edit script for Arbeit 7.docx
# Implements the diff: rewrites the "Problemstellung" closing paragraph,
# removes two paragraphs, turns the empty bold/tabs paragraph into the
# "Unsere Lösung für das Problem" heading (with list numbering), and
# appends two new paragraphs describing the solution.

$d = $word.ActiveDocument

function Find-ParagraphIndex($substr) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($substr)) {
            return $i
        }
    }
    return -1
}

# Common run properties used throughout this section of the document.
$rprPlain = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rprBold  = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

function New-PkgXml($bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Step 1: paragraph "Durch dieses Projekt möchten wir ... beitragen."
# Insert six new runs before the four existing ones, then delete the
# old run text so only the new runs remain.
# ---------------------------------------------------------------------
$idx1 = Find-ParagraphIndex("Durch dieses Projekt m")
$p1 = $d.Paragraphs($idx1)
$insertPoint = $d.Range($p1.Range.Start, $p1.Range.Start)

$newRuns1 = ''
$newRuns1 += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">Dabei ist das Transportieren schwerer Gegenstände, wie zum Beispiel Einkäufen, oftmals ein Problem von großer Bedeutung. </w:t></w:r>'
$newRuns1 += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">Dadurch können vor allem alleinlebende Menschen sich nicht mehr selbst versorgen und sind auf Hilfe durch Pflegekräfte angewiesen, welche viel Geld kosten, </w:t></w:r>'
$newRuns1 += '<w:r>' + $rprPlain + '<w:t>das</w:t></w:r>'
$newRuns1 += '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> man </w:t></w:r>'
$newRuns1 += '<w:r>' + $rprPlain + '<w:t>selten</w:t></w:r>'
$newRuns1 += '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> aufbringen kann.</w:t></w:r>'

# Wrapping in a bare <w:p> (no pPr) merges these runs into the existing
# target paragraph instead of creating a new paragraph break, because the
# insertion point is collapsed *inside* that paragraph.
$insertPoint.InsertXML((New-PkgXml ('<w:p>' + $newRuns1 + '</w:p>')))

# Re-locate the old run text (now shifted right by the inserted runs)
# and remove it, leaving only the six new runs behind.
$oldText1 = "Durch dieses Projekt m" + [char]0x00F6 + "chten wir auf das schwerwiegende Problem aufmerksam machen und unseren eigenen Teil zu der L" + [char]0x00F6 + "sung dieser Problematik beitragen."
$full = $d.Range(0, $d.Content.End).Text
$startOld = $full.IndexOf($oldText1)
$oldRange1 = $d.Range($startOld, $startOld + $oldText1.Length)
$oldRange1.Delete()

# ---------------------------------------------------------------------
# Step 2: paragraph "Unsere Idee ist es, ... abzunehmen." becomes empty
# (its pPr is preserved, only the run is removed).
# ---------------------------------------------------------------------
$idx2 = Find-ParagraphIndex("Unsere Idee ist es")
$p2 = $d.Paragraphs($idx2)
$ideaText = "Unsere Idee ist es, gerade diesen Menschen das Problem des Transportierens von beispielsweise Lebensmitteln nach einem Einkauf abzunehmen. "
$delRange2 = $d.Range($p2.Range.Start, $p2.Range.Start + $ideaText.Length)
$delRange2.Delete()

# ---------------------------------------------------------------------
# Step 3: paragraph "Durch eine möglichst simple ... auskennen." is
# removed entirely (paragraph mark included).
# ---------------------------------------------------------------------
$idx3 = Find-ParagraphIndex("Durch eine m")
$p3 = $d.Paragraphs($idx3)
$p3.Range.Delete()

# ---------------------------------------------------------------------
# Step 4: the empty bold/tabs paragraph becomes the new heading
# "Unsere Lösung für das Problem" with list numbering.
# ---------------------------------------------------------------------
$idx4 = Find-ParagraphIndex("")
# Locate the specific empty paragraph that carries the bold tab stop
# formatting (directly before the section break).
$n = $d.Paragraphs.Count
$idx4 = $n
$p4 = $d.Paragraphs($idx4)

$headingPPr = '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="9"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="1890"/></w:tabs>' + $rprBold + '</w:pPr>'
$headingPara = '<w:p>' + $headingPPr + '<w:r>' + $rprBold + '<w:t>Unsere L' + [char]0x00F6 + 'sung f' + [char]0x00FC + 'r das Problem</w:t></w:r></w:p>'

$p4.Range.InsertXML((New-PkgXml $headingPara))

$newHeadingIdx = $idx4
$oldEmptyIdx = $idx4 + 1
$newHeading = $d.Paragraphs($newHeadingIdx)
$oldEmpty = $d.Paragraphs($oldEmptyIdx)
$mergeRange1 = $d.Range($newHeading.Range.End - 1, $oldEmpty.Range.End)
$mergeRange1.Delete()

# ---------------------------------------------------------------------
# Step 5: append two new paragraphs after the heading: an empty spacing
# paragraph, then the paragraph describing the robot solution.
# ---------------------------------------------------------------------
$emptyPPr = '<w:pPr><w:spacing w:line="0" w:lineRule="atLeast"/>' + $rprPlain + '</w:pPr>'
$emptyPara = '<w:p>' + $emptyPPr + '</w:p>'

$solutionRuns = ''
$solutionRuns += '<w:r>' + $rprPlain + '<w:t>Unsere L' + [char]0x00F6 + 'sung</w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> daf' + [char]0x00FC + 'r</w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t xml:space="preserve"> besteht darin, einen autonom fahrenden Roboter </w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">zu entwickeln, welcher eine </w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">die Eink' + [char]0x00E4 + 'ufe f' + [char]0x00FC + 'r einen </w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t>transportiert</w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t xml:space="preserve">, </w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t>sodass man selbst N</w:t></w:r>'
$solutionRuns += '<w:r>' + $rprPlain + '<w:t>ichts schweres mehr tragen muss.</w:t></w:r>'
$solutionPara = '<w:p>' + $emptyPPr + $solutionRuns + '</w:p>'

$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIdx)
$pLast.Range.InsertXML((New-PkgXml ($emptyPara + $solutionPara)))

$newSolutionIdx = $lastIdx + 1
$oldEmptyIdx2 = $lastIdx + 2
$newSolution = $d.Paragraphs($newSolutionIdx)
$oldEmpty2 = $d.Paragraphs($oldEmptyIdx2)
$mergeRange2 = $d.Range($newSolution.Range.End - 1, $oldEmpty2.Range.End)
$mergeRange2.Delete()

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
